$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "general" sheet: insert a new individual (upf_00205 / "FJSL young brother")
# as a new row 6, pushing the existing rows 6-11 down to 7-12.
# ---------------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")

# Drop the stray "married_with" value that used to sit on FJSL Spouse's row
# (D5) - it belongs on John Doe's row after the insert below.
$general.Range("D5").ClearContents()

$general.Rows.Item(6).Insert(-4121)
$general.Range("A6").Value = "upf_00205"
$general.Range("B6").Value = "FJSL young brother"
$general.Range("C6").Value = "male"
$general.Range("D6").Value = "upf_00201"
$general.Range("E6").Value = "upf_00202"
$general.Range("G6").Value = "upf_f_2"
$general.Range("H6").Value = "UPF"

# ---------------------------------------------------------------------------
# "clinical values" sheet: mirror the same insertion with its clinical data.
# ---------------------------------------------------------------------------
$clinical = $wb.Worksheets.Item("clinical values")

$clinical.Rows.Item(6).Insert(-4121)
$clinical.Range("A6").Value = "upf_00205"
$clinical.Range("C7").Copy()
$clinical.Range("C6").PasteSpecial(-4122)
$clinical.Range("C6").Value = 28856
$clinical.Range("E2").Copy()
$clinical.Range("E6").PasteSpecial(-4163)
$clinical.Range("F2").Copy()
$clinical.Range("F6").PasteSpecial(-4163)
$clinical.Range("I2").Copy()
$clinical.Range("I6").PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# Viewer sidebar state: the "general" sheet becomes the active tab, with the
# cursor left on E8; "Family" is no longer the selected tab.
# ---------------------------------------------------------------------------
$general.Select()
$general.Range("E8").Select()
$clinical.Range("K7").Select()
$general.Select()
